# "fixed mistake in layout"
# The ConAssembly stations' ABU4 utilisation figures were wrong: ConAssembly1 (row 11)
# was missing its ABU4 load and ConAssembly2 (row 12) had it overstated. Both should be 0.5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E11").Value = 0.5
$ws.Range("E12").Value = 0.5

# The "Description" column (H) was far wider than the text needs, squeeze it back down.
$ws.Columns.Item(8).ColumnWidth = 102.67

# Leave the cursor where the author ended up after making the fix.
$ws.Range("F22").Select()
